$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Trip Type" column values from "one way" to "One Way"
$ws.Range("D2:D5").Value = "One Way"

# Reflect the active cell selection recorded in the saved workbook
$ws.Range("I3").Select()
